$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select range B3:F5 and clear its contents (mirrors a user selecting the
# range and pressing Delete), which leaves the A column labels intact while
# removing the description/time/location/link cells in rows 3-5.
$range = $ws.Range("B3:F5")
$range.Select()
$range.ClearContents()
